{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is exactly \"\u661f\u671f\u5929\" and the one whose\n// text is exactly \"\u661f\u671f\u4e00\" (the last paragraph in the original document).\nlet weekSunday = null;\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"\u661f\u671f\u5929\") {\n    weekSunday = paragraphs.items[i];\n  } else if (t === \"\u661f\u671f\u4e00\") {\n    target = paragraphs.items[i];\n  }\n}\n\nif (!weekSunday) {\n  throw new Error(\"Could not find paragraph with text \u661f\u671f\u5929\");\n}\nif (!target) {\n  throw new Error(\"Could not find paragraph with text \u661f\u671f\u4e00\");\n}\n\n// Insert a fresh paragraph reading \"\u661f\u671f\u4e00\" right after \"\u661f\u671f\u5929\" (it picks\n// up that paragraph's plain formatting, with no extra pPr), so the\n// document now reads ... \u661f\u671f\u5929 / \u661f\u671f\u4e00(new) / \u661f\u671f\u4e00(old).\nweekSunday.insertParagraph(\"\u661f\u671f\u4e00\", \"After\");\n\n// Turn the original \"\u661f\u671f\u4e00\" paragraph into the new closing line.\ntarget.insertText(\"\u4eca\u5929\u5929\u6c14\u4e0d\u9519\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph reading \"\u661f\u671f\u5929\" and the one reading \"\u661f\u671f\u4e00\" (the\n# last paragraph in the original document). Paragraph.Range.Text includes\n# the trailing paragraph mark, so compare against that form.\n$weekSunday = $null\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -eq \"\u661f\u671f\u5929`r\") {\n        $weekSunday = $p\n    }\n    if ($t -eq \"\u661f\u671f\u4e00`r\") {\n        $target = $p\n    }\n}\n\nif ($weekSunday -eq $null) {\n    throw \"Could not find paragraph with text \u661f\u671f\u5929\"\n}\nif ($target -eq $null) {\n    throw \"Could not find paragraph with text \u661f\u671f\u4e00\"\n}\n\n# Turn the original \"\u661f\u671f\u4e00\" paragraph into the new closing line first,\n# before the paragraph collection shifts because of the insert below.\n$target.Range.Text = \"\u4eca\u5929\u5929\u6c14\u4e0d\u9519\"\n\n# Insert a fresh paragraph right after \"\u661f\u671f\u5929\" (it picks up that\n# paragraph's plain formatting, with no extra pPr) and give it the text\n# \"\u661f\u671f\u4e00\", so the document now reads ... \u661f\u671f\u5929 / \u661f\u671f\u4e00(new) / \u4eca\u5929\u5929\u6c14\u4e0d\u9519.\n$weekSunday.Range.InsertParagraphAfter()\n$weekSunday.Next().Range.Text = \"\u661f\u671f\u4e00\"\n\n$d.Save()\n"}
